# Bug and Task List - update footer in project
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Status (column D) updates -------------------------------------------
$ws.Range("D2").Value  = "In-Progress"
$ws.Range("D3").Value  = "Done"
$ws.Range("D4").Value  = "Done"
$ws.Range("D5").Value  = "Done"
$ws.Range("D13").Value = "Done"
$ws.Range("D14").Value = "In-Progress"

# --- Remove the AutoFilter and reveal the rows it had hidden --------------
$ws.Rows.Item(6).Hidden  = $false
$ws.Rows.Item(7).Hidden  = $false
$ws.Rows.Item(8).Hidden  = $false
$ws.Rows.Item(12).Hidden = $false
$ws.Rows.Item(15).Hidden = $false

if ($ws.AutoFilterMode) {
    $ws.AutoFilterMode = $false
}

# --- Append two new task rows ---------------------------------------------
$ws.Range("B16").Copy() | Out-Null
$ws.Range("B17").PasteSpecial(-4122) | Out-Null
$ws.Range("B17").Value = "Changes in tax validation and implementation"

$ws.Range("B16").Copy() | Out-Null
$ws.Range("B18").PasteSpecial(-4122) | Out-Null
$ws.Range("B18").Value = "CP1 Quality development"

$ws.Range("A17").Value = 16
$ws.Range("C17").Value = "Code"
$ws.Range("D17").Value = "Pending"

$ws.Range("A18").Value = 17
$ws.Range("C18").Value = "Code"
$ws.Range("D18").Value = "Pending"

$ws.Range("A19").Value = 18

# --- Selection matches the authored file -----------------------------------
$ws.Range("D18").Select() | Out-Null

Write-Host "Applied edits"
